$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Replace the text of the final (italic) paragraph with the new
# "Create a feature image ..." image-prompt copy, while leaving its run
# formatting (the <w:i/> run) and its leading empty run untouched.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$italicPara = $d.Paragraphs.Item($count)

$newImagePrompt = "Create a feature image that captures the all-powerful and mischievous Baron Samedi in a cartoon-style design. The image should include a happy Maya warrior wearing glasses, as this character represents the player in the game. Be sure to include elements of voodoo, such as candles, bones, and alcohol, to set the scene. The image should be bold and eye-catching, with vibrant colors that reflect the mystical and mysterious world of voodoo. It should draw players in and make them curious to try out the game and experience the fun and excitement of playing with the charismatic Baron Samedi."

$pStart = $italicPara.Range.Start
$pEnd = $italicPara.Range.End - 1
$textRange = $d.Range($pStart, $pEnd)
$textRange.Text = $newImagePrompt

# ---------------------------------------------------------------------------
# Step 2: Delete the old bold "Play Baron Samedi Free - Exciting Voodoo
# Themed Slot" paragraph that used to sit right before the paragraph handled
# above (it is being dropped outright, not just its text).
# ---------------------------------------------------------------------------
$boldPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$boldPara.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# Step 3: Insert a brand-new paragraph right after the document's opening
# Heading1 paragraph, containing a bold "Meta description" run followed by
# a plain run with the rest of the description text.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphAfter() | Out-Null
$metaPara = $d.Paragraphs.Item(2)
$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Baron Samedi - the latest addition to the voodoo-themed slot games. Play Baron Samedi free and enjoy the unique gaming experience it provides.</w:t></w:r></w:p>'
$metaPara.Range.InsertXML($metaXml) | Out-Null
